# Edit script for 016.content.docx
# Implements:
#  1. Remove the "License Information" Heading2 paragraph.
#  2. Replace the license paragraph's text with the new Biblica Study Notes
#     copyright/credit text, and remove the following
#     "This PDF version is provided under the same license." paragraph.
#  3. Remove the italic "key terms" index paragraph that used to follow the
#     "ك" Heading2 paragraph.

$d = $word.ActiveDocument

function Find-ParagraphByText($doc, [string]$needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        if ($para.Range.Text.Contains($needle)) {
            return $para
        }
    }
    return $null
}

# --- Step 1: delete the "License Information" paragraph entirely ---
$licPara = Find-ParagraphByText $d "License Information"
if ($licPara -ne $null) {
    $licPara.Range.Delete()
}

# --- Step 2: rewrite the license-credit paragraph ---
$creditPara = Find-ParagraphByText $d "is based on"
$pStart = $creditPara.Range.Start
$pEnd = $creditPara.Range.End

# Replace the bold heading run's text.
$boldRng = $d.Range($pStart, $pEnd)
$boldRng.Find.Execute("المصطلحات الرئيسية (Biblica)", $false, $false, $false, $false, $false, $true, 1, $false, "Biblica Study Notes (Key Terms)", 2) | Out-Null
$afterBold = $boldRng.End

# Clear everything after the bold run through the end of the paragraph text
# (but not the paragraph mark itself, which is the very last character).
$pEndNow = $creditPara.Range.End
$tailRng = $d.Range($afterBold, $pEndNow - 1)
$tailRng.Text = ""

# Insert the new tail text as a single non-bold run.
$newTail = " © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. Biblica Study Notes has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual."
$insRng = $d.Range($afterBold, $afterBold)
$insRng.Text = $newTail
$insRng.Font.Bold = 0

# --- Step 2b: delete the following "This PDF version..." paragraph ---
$pdfPara = Find-ParagraphByText $d "This PDF version is provided under the same license."
if ($pdfPara -ne $null) {
    $pdfPara.Range.Delete()
}

# --- Step 3: delete the italic key-terms index paragraph ---
$listPara = Find-ParagraphByText $d "الكائنات الروحية الشريرة, كائنات روحية, كالب"
if ($listPara -ne $null) {
    $listPara.Range.Delete()
}
